$d = $word.ActiveDocument

# 1. Remove the red highlight from the "Estándar: ..." run of text.
#    Setting HighlightColorIndex on any sub-range clears it for all the
#    (highlighted) runs sharing that paragraph's text.
$rng = $d.Content
$found = $rng.Find.Execute("Identiﬁco", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.HighlightColorIndex = 0
}

# 2. Drop the trailing period from "Relaciones con la historia y las culturas."
$rng = $d.Content
$rng.Find.Execute("Relaciones con la historia y las culturas.", $true, $false, $false, $false, $false, $true, 1, $false, "Relaciones con la historia y las culturas", 2) | Out-Null

# 3. Split "Relaciones" into its own run (without changing formatting).
$rng = $d.Content
$found = $rng.Find.Execute("Relaciones", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Font.Bold = 0
    $rng.Font.Bold = 1
}

# 4. Split "histo" / "ria" and move the _GoBack bookmark to sit between them.
$rng = $d.Content
$found = $rng.Find.Execute("histo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $rng)
}

# 5. Re-touch the trailing "ria y las culturas" run so it is re-serialized
#    cleanly (avoids a spurious xml:space="preserve").
$rng = $d.Content
$rng.Find.Execute("ria y las culturas", $true, $false, $false, $false, $false, $true, 1, $false, "ria y las culturas", 2) | Out-Null
